# Apply odds updates to Sheet1 as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5
$ws.Range("Q5").Value = 1.5

# Row 8
$ws.Range("F8").Value = 2.64
$ws.Range("G8").Value = 2.84
$ws.Range("Q8").Value = 1.68

# Row 9
$ws.Range("F9").Value = 2.06
$ws.Range("H9").Value = 3.9
$ws.Range("I9").Value = 4.8
$ws.Range("Q9").Value = 1.72

# Row 10
$ws.Range("H10").Value = 1.74
$ws.Range("I10").Value = 1.75
$ws.Range("O10").Value = 1.31
$ws.Range("Q10").Value = 1.92
$ws.Range("R10").Value = 1.41
$ws.Range("S10").Value = 3.35
$ws.Range("T10").Value = 1.89
$ws.Range("X10").Value = 15.5
$ws.Range("AA10").Value = 18
$ws.Range("AJ10").Value = 140
$ws.Range("AM10").Value = 120
$ws.Range("AO10").Value = 11

# Row 11
$ws.Range("F11").Value = 1.91
$ws.Range("G11").Value = 2.14
$ws.Range("H11").Value = 4.6
$ws.Range("J11").Value = 3.1
$ws.Range("K11").Value = 3.55
$ws.Range("P11").Value = 1.53
$ws.Range("Q11").Value = 2.54

# Row 12
$ws.Range("F12").Value = 2.26
$ws.Range("H12").Value = 3.65
$ws.Range("I12").Value = 4
$ws.Range("J12").Value = 3.15
$ws.Range("P12").Value = 1.71
$ws.Range("Q12").Value = 2.2

# Row 13
$ws.Range("H13").Value = 3.8
$ws.Range("K13").Value = 3.45
$ws.Range("Q13").Value = 2.12

# Row 14
$ws.Range("G14").Value = 2.24

# Row 16
$ws.Range("F16").Value = 2.22
$ws.Range("G16").Value = 2.46
$ws.Range("I16").Value = 3.85
